$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin (D2 already non-numeric-looking text, keep as-is)
$ws.Range("D2").Value = "64.329.88"
$ws.Range("E2").Value = "  +1.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.672.18"
$ws.Range("E3").Value = "  +2.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB (numeric-looking -> force text with leading apostrophe)
$ws.Range("D5").Value = "'596.07"
$ws.Range("E5").Value = "  +2.09%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'147.68"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.99%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  +0.13%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -0.28%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.03%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +0.65%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'27.87"
$ws.Range("E13").Value = "  +2.23%  "

# Row 14 - WrappedliquidstakedEther2.0 (multi-dot, naturally text)
$ws.Range("D14").Value = "3.154.56"
$ws.Range("E14").Value = "  +3.06%  "

# Row 15 - WrappedBTC (multi-dot, naturally text)
$ws.Range("D15").Value = "64.246.44"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.14%  "

# Row 17 - WrappedEther (multi-dot, naturally text)
$ws.Range("D17").Value = "2.723.46"
$ws.Range("E17").Value = "  +5.29%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'11.42"
$ws.Range("E18").Value = "  +0.44%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'346.70"
$ws.Range("E19").Value = "  +0.91%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'4.40"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "  +1.33%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.16%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'68.75"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24 - SuiNetwork
$ws.Range("D24").Value = "'1.61"
$ws.Range("E24").Value = "  +9.94%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  +4.04%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -1.41%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'8.53"
$ws.Range("E27").Value = "  +1.64%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "'7.99"
$ws.Range("E28").Value = "  +1.10%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.23%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "'529.90"
$ws.Range("E30").Value = "  +13.83%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'2.00"
$ws.Range("E31").Value = "  +3.08%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.78"
$ws.Range("E32").Value = "  +11.13%  "

# Row 33 - PEPE (contains unicode subscript, naturally text)
$ws.Range("D33").Value = "0.0₃0827"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34 - Monero
$ws.Range("D34").Value = "'175.52"
$ws.Range("E34").Value = "  -0.74%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.07%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("D36").Value = "'0.403"

# Row 37 - EthereumClassic
$ws.Range("E37").Value = "  +0.63%  "

# Row 38 - NEARProtocol
$ws.Range("D38").Value = "'4.70"
$ws.Range("E38").Value = "  +2.06%  "

# Row 39 - was Stacks, now Aave
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "'173.27"
$ws.Range("E39").Value = "  +8.97%  "

# Row 40 - was Aave, now Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "  +3.11%  "

# Row 41 - USDe
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - OKB
$ws.Range("D42").Value = "'40.72"
$ws.Range("E42").Value = "  +3.11%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  -0.38%  "

# Row 44 - InjectiveProtocol
$ws.Range("D44").Value = "'21.74"
$ws.Range("E44").Value = "  +2.79%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "'0.634"
$ws.Range("E45").Value = "  -0.50%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0550"
$ws.Range("E46").Value = "  +0.33%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "'0.0241"
$ws.Range("E47").Value = "  +1.24%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -1.10%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'18.86"
$ws.Range("E49").Value = "  +1.52%  "

# Row 50 - dogwifhat
$ws.Range("E50").Value = "  +2.52%  "

# Row 51 - WhiteBITCoin
$ws.Range("D51").Value = "'11.35"
$ws.Range("E51").Value = "  -0.67%  "
